# Revert "Adding Passing Score to Assessment, excel change STRING to text"
# -> change the TextType column (E2:E17) back from "text" to "STRING"
# and leave the selection spanning the full used range (A1:P17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E17").Value = "STRING"

$ws.Range("A1:P17").Select()
